# Adapt the worksheet to the new file format: add a "Grupa (prowadzący)"
# question block (3 columns: answer / points / feedback) right after the
# existing "Nazwa zwierzaka" question block and before the
# "Ustaw nagrody..." question block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns at Q:S. This pushes the existing Q:S
# ("Ustaw nagrody wg preferowanej kolejności..." block) to T:V.
$ws.Columns("Q:S").Insert()

# New header row for the inserted "Grupa (prowadzący)" question block.
$ws.Cells.Item(1, 17).Value = "Grupa (prowadzący):"
$ws.Cells.Item(1, 18).Value = "Punkty — Grupa (prowadzący):"
$ws.Cells.Item(1, 19).Value = "Opinia — Grupa (prowadzący):"

# Answers for each respondent row.
$ws.Cells.Item(2, 17).Value = "Bernard Maj"
$ws.Cells.Item(3, 17).Value = "Zbigniew Kaleta"
$ws.Cells.Item(4, 17).Value = "Michał Idzik"

# Match the final selection left in the saved file.
$ws.Range("T9").Select()
